# Rename/reorder primary-key label cells in the StarSchemaPlan sheet.
# (commit: "Added a copy of Airflow dag")
#
# Changes (column B only; A/C/D/E on these rows are untouched):
#   B42 : ProviderID(PK)            -> PhysicianID(PK)
#   B61 : DiagnosisCode              -> DiagnosisCodeID(PK)
#   B62 : DiagnosisCodeDescription  -> DiagnosisCodeGroup
#   B63 : DiagnosisCodeGroup         -> DiagnosisCodeDescription
#   B65 : CPTCode                   -> CPTCodeID(PK)
#   B66 : CPTCodeDescription        -> CPTCodeGroup
#   B67 : CPTCodeGroup               -> CPTCodeDescription

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B42").Value = "PhysicianID(PK)"

$ws.Range("B61").Value = "DiagnosisCodeID(PK)"
$ws.Range("B62").Value = "DiagnosisCodeGroup"
$ws.Range("B63").Value = "DiagnosisCodeDescription"

$ws.Range("B65").Value = "CPTCodeID(PK)"
$ws.Range("B66").Value = "CPTCodeGroup"
$ws.Range("B67").Value = "CPTCodeDescription"

# Match the saved selection state from the edit.
$ws.Activate()
$ws.Range("B73").Select() | Out-Null
